# Added config for Passing base URL as Parameter
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Part Of Test Plan" (column D) values that flipped.
$ws.Range("D2").Value = "N"
$ws.Range("D20").Value = "Y"
$ws.Range("D21").Value = "N"
$ws.Range("D22").Value = "N"

# Move the active selection to D20 to match the saved view state.
$ws.Range("D20").Select()
